$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.968.21'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '3.112.23'
$ws.Range("E3").Value = '  +2.53%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'579.41"
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").Value = "'172.95"
$ws.Range("E6").Value = '  +2.67%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.106.14'
$ws.Range("E8").Value = '  +2.51%  '
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("E10").Value = '  -3.11%  '
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").Value = "'37.29"
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("D15").Value = '0.124'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '3.627.86'
$ws.Range("E16").Value = '  +2.26%  '
$ws.Range("D17").Value = '67.016.41'
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").Value = "'7.20"
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '3.112.30'
$ws.Range("E19").Value = '  +2.16%  '
$ws.Range("D20").Value = "'16.30"
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Value = "'486.55"
$ws.Range("E21").Value = '  +4.29%  '
$ws.Range("D22").Value = "'0.720"
$ws.Range("E22").Value = '  +2.12%  '
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("D24").Value = "'84.61"
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("D25").Value = "'13.36"
$ws.Range("E25").Value = '  +3.56%  '
$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = '  +4.07%  '
$ws.Range("D27").Value = '10.09'
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").Value = "'8.07"
$ws.Range("E29").Value = '  -3.50%  '
$ws.Range("D30").Value = "'2.43"
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("E31").Value = '  +2.49%  '
$ws.Range("D32").Value = "'29.01"
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("E34").Value = '  -3.13%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").Value = "'5.95"
$ws.Range("E36").Value = '  +2.17%  '
$ws.Range("D37").Value = "'0.988"
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("B38").Value = 'Arweave'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D38").Value = "'47.41"
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = "'2.12"
$ws.Range("E39").Value = '  +3.84%  '
$ws.Range("D40").Value = "'50.13"
$ws.Range("E40").Value = '  +1.26%  '
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("D43").Value = "'8.69"
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("D45").Value = '2.843.31'
$ws.Range("E45").Value = '  +3.94%  '
$ws.Range("D46").Value = '385.44'
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = "'0.0359"
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").Value = '136.37'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = "'25.08"
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("E51").Value = '  +0.45%  '
